# PWR_Board_TestReportTemplate2.xlsx - "Report" sheet updates
# Commit: re-implemented get_current_protection_state / fixed HV diode tests,
# clearing out the -999 MEASURED placeholders and updating the diode
# threshold rows (10-15) to literal pass/fail values instead of the
# shared "1.1 * F" formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")
$ws.Activate()

# The "MEASURED" column (G) was pre-seeded with a -999 sentinel on every
# test row; clear those placeholders back out to empty cells (rows 2-95).
$ws.Range("G2:G95").ClearContents()

# Diode threshold rows (9 -> rows 10-15): drop the shared formula
# "=1.1*F10" (which depended on the old 0.05 MAX column) and replace both
# the MAX (E) and NOMINAL (F) columns with their fixed literal values.
$ws.Range("E10:E15").Value = 0.001
$ws.Range("F10:F15").Value = 0

# Leave the selection where the author was last working.
$ws.Range("E10:E15").Select()
